$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 6)
$ws.Range("M6").Value = "precision"
$ws.Range("O6").Value = "log-based precision"
$ws.Range("R6").Value = "fitness"

# Data rows 7-24: M = precision value, O = formula M*2, R = fitness value
$data = [ordered]@{
    7  = @(0.4347, 0.9436)
    8  = @(0.4,    0.9542)
    9  = @(0.4347, 0.9448)
    10 = @(0.25,   0.9921)
    11 = @(0.4,    0.9278)
    12 = @(0.3,    0.9463)
    13 = @(0.4,    1)
    14 = @(0.413,  0.9272)
    15 = @(0.3,    0.9982)
    16 = @(0.4782, 0.953)
    17 = @(0,      1)
    18 = @(0,      1)
    19 = @(0.4782, 0.9565)
    20 = @(0,      1)
    21 = @(0.4782, 0.9475)
    22 = @(0,      1)
    23 = @(0,      1)
    24 = @(0,      1)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Range("M$r").Value = $vals[0]
    $ws.Range("O$r").Formula = "=M$r*2"
    $ws.Range("R$r").Value = $vals[1]
}

$ws.Range("O15:O16").Select()
